$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet Paquete -> Precios
$ws.Name = "Precios"

# 2. Defined names: Estudios now points at the moved table; Paquete -> Precios
#    (Precios local name refers to an external workbook copy of this sheet)
$wb.Names.Item("Estudios").Delete()
$wb.Names.Item("Precios!Paquete").Delete()
$wb.Names.Add("Estudios", "=Precios!`$B`$15:`$E`$16")
$ws.Names.Add("Precios", "=[1]Precios!`$A`$3:`$C`$16")

# 3. Update placeholder text (Paquete.* -> Precios.*) for the surviving fields
$ws.Range("B3").Value = "{{Precios.Clave}}"
$ws.Range("B5").Value = "{{Precios.Nombre}}"

# Row 7 used to be "Nombre largo" -> now "Activo"
$ws.Range("A7").Value = "Activo"
$ws.Range("B7").Value = "{{Precios.Activo}}"

# Row 9 used to be "Área" -> now "Visible"
$ws.Range("A9").Value = "Visible"
$ws.Range("B9").Value = "{{Precios.visible}}"

# Row 11 used to be "Departamento" -> now blank (labels removed)
$ws.Range("A11").Value = ""
$ws.Range("B11").Value = ""

# 4. Remove the old "Estudios" mini-table block (rows 20-23) entirely
$ws.Range("A20:F23").EntireRow.Delete()

# 5. Rebuild rows 13-16 as the relocated/extended Estudios table
$ws.Range("A13:C13").UnMerge()
$ws.Range("A13").Clear()
$ws.Range("B13:C15").UnMerge()

$ws.Range("B13:D13").Merge()
$ws.Range("B13").Value = "Estudios"

$ws.Range("B14").Value = "Clave"
$ws.Range("C14").Value = "Nombre"
$ws.Range("D14").Value = "Área"
$ws.Range("E14").Value = "Precio"

$ws.Range("B15").Value = "{{item.Clave}}"
$ws.Range("C15").Value = "{{item.Nombre}}"
$ws.Range("D15").Value = "{{item.Area}}"
$ws.Range("E15").Value = "{{item.Precio}}"

# 6. Selection moves to D3, no frozen/scrolled topLeftCell anymore
$ws.Range("D3").Select() | Out-Null
